$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I ("sexo") is re-curated from an sdmx/iaest DIMENSION to a MEASURE:
#   iaest-dimension:sexo -> iaest-measure:sexo
#   dim                  -> medida
#   skos:Concept         -> xsd:int
# and its mapping file reference (row 5) is removed entirely (cell dropped).
$ws.Range("I2").Value = "iaest-measure:sexo"
$ws.Range("I3").Value = "medida"
$ws.Range("I4").Value = "xsd:int"
$ws.Range("I5").Clear()

# Column K ("direccion-provincial-nombre") is re-curated from a generic
# refArea dimension to its own measure:
#   sdmx-dimension:refArea -> iaest-measure:direccion-provincial-nombre
#   dim                     -> medida
#   URI-Provincia           -> xsd:int
$ws.Range("K2").Value = "iaest-measure:direccion-provincial-nombre"
$ws.Range("K3").Value = "medida"
$ws.Range("K4").Value = "xsd:int"
